$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 235
$ws.Range("B3").Value = 320
$ws.Range("B4").Value = 150
$ws.Range("B5").Value = 200
